$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 63: last name entered before first name (matches original authoring order)
$ws.Cells.Item(63, 2).Value = "Merhbene"
$ws.Cells.Item(63, 1).Value = "Ghofrane"
$ws.Cells.Item(63, 3).Value = "Bern University of Applied Sciences"
$ws.Cells.Item(63, 4).Value = "Suisse"
$ws.Cells.Item(63, 5).Value = "I5KRVLMAAAAJ"
$ws.Cells.Item(63, 6).Value = "F"
$ws.Cells.Item(63, 7).Value = 1999
$ws.Cells.Item(63, 8).Value = "Informatique, Mathématiques et Ingénierie"

# Row 64: normal first-name, last-name order
$ws.Cells.Item(64, 1).Value = "Amr"
$ws.Cells.Item(64, 2).Value = "Chaabani"
$ws.Cells.Item(64, 3).Value = "Université de Monastir"
$ws.Cells.Item(64, 4).Value = "Tunisie"
$ws.Cells.Item(64, 5).Value = "eickZDMAAAAJ"
$ws.Cells.Item(64, 6).Value = "M"
$ws.Cells.Item(64, 7).Value = 1994
$ws.Cells.Item(64, 8).Value = "Médecine, Biologie et Sciences de la Santé"

# Row 65: normal first-name, last-name order
$ws.Cells.Item(65, 1).Value = "Dorra"
$ws.Cells.Item(65, 2).Value = "Hadj Mahmoud"
$ws.Cells.Item(65, 3).Value = "Université d'Angers"
$ws.Cells.Item(65, 4).Value = "France"
$ws.Cells.Item(65, 5).Value = "20soGN4AAAAJ"
$ws.Cells.Item(65, 6).Value = "F"
$ws.Cells.Item(65, 7).Value = 1990
$ws.Cells.Item(65, 8).Value = "Médecine, Biologie et Sciences de la Santé"

# Apply the same style used for the other Genre (column F) cells
$ws.Range("F62").Copy()
$ws.Range("F63:F65").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Application.ActiveWindow.ScrollRow = 39
$ws.Range("H66").Select()
